$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Pass" result to column G for rows 7 and 8, matching the bordered
# style already used by column A on those rows (xf with left/right borders).
$ws.Range("A7").Copy()
$ws.Range("G7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G7").Value = "Pass"

$ws.Range("A8").Copy()
$ws.Range("G8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G8").Value = "Pass"

$excel.CutCopyMode = $false

# Update the active selection on the sheet view
$ws.Range("A2:G5").Select()
